$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-23"

# Update the label in A7 ("June (through 06-22)" -> "June (through 06-23)")
$ws.Range("A7").Value = "June (through 06-23)"

# Update the June row (row 7) values
$ws.Range("D7").Value = 56
$ws.Range("E7").Value = 44
$ws.Range("G7").Value = 90
$ws.Range("H7").Value = 89
$ws.Range("I7").Value = 108

# Update the Total row (row 8) values
$ws.Range("D8").Value = 372
$ws.Range("E8").Value = 339
$ws.Range("G8").Value = 448
$ws.Range("H8").Value = 720
$ws.Range("I8").Value = 771
